$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a floating point re-computation of the existing A6 timestamp
$ws.Cells.Item(6, 1).Value = 44319.77475232986

# Append new row 7 with data retrieved on 2021-05-04 (see commit message)
$ws.Cells.Item(7, 1).Value = 44320.77039254122
$ws.Cells.Item(7, 1).NumberFormat = $ws.Cells.Item(6, 1).NumberFormat

$ws.Cells.Item(7, 2).Value = 71249
$ws.Cells.Item(7, 3).Value = 59997
$ws.Cells.Item(7, 4).Value = 3142
$ws.Cells.Item(7, 5).Value = 1967
$ws.Cells.Item(7, 6).Value = 1398
$ws.Cells.Item(7, 7).Value = 18719
$ws.Cells.Item(7, 8).Value = 1298
$ws.Cells.Item(7, 9).Value = 808
$ws.Cells.Item(7, 10).Value = 193
